# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "40.164.48"
$ws.Range("E2").Value = "  +0.09%  "
$ws.Range("D3").Value = "2.212.72"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'295.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.09%  "
$ws.Range("D6").Value = "'87.73"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.19%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.71%  "
$ws.Range("D10").Value = "'30.92"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.51%  "
$ws.Range("D11").Value = "'51.47"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +6.00%  "
$ws.Range("D12").Value = "'0.0781"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.08%  "
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "2.553.29"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").Value = "'13.82"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.14%  "
$ws.Range("D17").Value = "2.208.11"
$ws.Range("E17").Value = "  -1.66%  "
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "40.066.52"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "0.0₃0887"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "'11.30"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("D23").Value = "'65.58"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").Value = "'235.52"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.62%  "
$ws.Range("D25").Value = "'0.999"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("D26").Value = "'2.49"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +0.90%  "
$ws.Range("E27").Value = "  -0.92%  "
$ws.Range("D28").Value = "'23.20"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.89%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.08"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.77%  "
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").Value = "'9.32"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").Value = "'157.10"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.37%  "
$ws.Range("D32").Value = "'32.08"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.74%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("E35").Value = "  +5.59%  "
$ws.Range("E36").Value = "  -0.50%  "
$ws.Range("E37").Value = "  -0.62%  "
$ws.Range("E38").Value = "  +1.32%  "
$ws.Range("E39").Value = "  +1.88%  "
$ws.Range("D40").Value = "'1.74"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.47%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").Value = "2.075.36"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("D43").Value = "'3.80"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.74%  "
$ws.Range("D44").Value = "'19.49"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +9.08%  "
$ws.Range("D45").Value = "'0.0272"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("D46").Value = "'9.94"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  +5.06%  "
$ws.Range("E48").Value = "  -10.62%  "
$ws.Range("D49").Value = "2.426.93"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'1.48"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.80%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.12"
$ws.Range("D51").ClearFormats()
